# edit.ps1
#
# Refreshes the BP Terminal Gate Pricing workbook for the new daily cycle:
# each terminal's two most-recent-date rows are shifted to the next effective
# date and the Diesel/ULP/PULP/e10 prices (columns D-G) are updated to the new
# published cents-per-litre values. Only data values change; formatting/styles
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8: Sydney-Botany
$ws.Range("A8").Value = 45953
$ws.Range("D8").Value = 159.62
$ws.Range("E8").Value = 157.68
$ws.Range("F8").Value = 167.68
$ws.Range("G8").Value = 157.83000000000001

# Row 9: Sydney-Silverwater
$ws.Range("A9").Value = 45953
$ws.Range("D9").Value = 159.62
$ws.Range("E9").Value = 157.68
$ws.Range("F9").Value = 167.68
$ws.Range("G9").Value = 157.83000000000001

# Row 10: Newcastle
$ws.Range("A10").Value = 45953
$ws.Range("D10").Value = 161.88999999999999
$ws.Range("E10").Value = 159.75
$ws.Range("F10").Value = 169.75
$ws.Range("G10").Value = 160.22

# Row 11: Sydney-Botany
$ws.Range("A11").Value = 45952
$ws.Range("D11").Value = 160.06
$ws.Range("E11").Value = 157.88999999999999
$ws.Range("F11").Value = 167.89
$ws.Range("G11").Value = 158.05000000000001

# Row 12: Sydney-Silverwater
$ws.Range("A12").Value = 45952
$ws.Range("D12").Value = 160.06
$ws.Range("E12").Value = 157.88999999999999
$ws.Range("F12").Value = 167.89
$ws.Range("G12").Value = 158.05000000000001

# Row 13: Newcastle
$ws.Range("A13").Value = 45952
$ws.Range("D13").Value = 162.28
$ws.Range("E13").Value = 160.22999999999999
$ws.Range("F13").Value = 170.23
$ws.Range("G13").Value = 160.69999999999999

# Row 17: Darwin
$ws.Range("A17").Value = 45953
$ws.Range("D17").Value = 165.32
$ws.Range("E17").Value = 162.80000000000001
$ws.Range("F17").Value = 172.8

# Row 18: Darwin
$ws.Range("A18").Value = 45952
$ws.Range("D18").Value = 165.7
$ws.Range("E18").Value = 163.07
$ws.Range("F18").Value = 173.07

# Row 22: Brisbane
$ws.Range("A22").Value = 45953
$ws.Range("D22").Value = 160.54
$ws.Range("E22").Value = 158.75
$ws.Range("F22").Value = 168.35
$ws.Range("G22").Value = 159.93

# Row 23: Cairns
$ws.Range("A23").Value = 45953
$ws.Range("D23").Value = 166.65
$ws.Range("E23").Value = 163.57
$ws.Range("F23").Value = 173.57

# Row 24: Gladstone
$ws.Range("A24").Value = 45953
$ws.Range("D24").Value = 166.46
$ws.Range("E24").Value = 163.78
$ws.Range("F24").Value = 173.78

# Row 25: Mackay
$ws.Range("A25").Value = 45953
$ws.Range("D25").Value = 167.29
$ws.Range("E25").Value = 163.16999999999999
$ws.Range("F25").Value = 173.17
$ws.Range("G25").Value = 163

# Row 26: Townsville
$ws.Range("A26").Value = 45953
$ws.Range("D26").Value = 166.01
$ws.Range("E26").Value = 164.71
$ws.Range("F26").Value = 174.71

# Row 27: Brisbane
$ws.Range("A27").Value = 45952
$ws.Range("D27").Value = 160.97999999999999
$ws.Range("E27").Value = 159.12
$ws.Range("F27").Value = 168.72
$ws.Range("G27").Value = 160.30000000000001

# Row 28: Cairns
$ws.Range("A28").Value = 45952
$ws.Range("D28").Value = 167.04
$ws.Range("E28").Value = 163.83000000000001
$ws.Range("F28").Value = 173.83

# Row 29: Gladstone
$ws.Range("A29").Value = 45952
$ws.Range("D29").Value = 166.85
$ws.Range("E29").Value = 164.04
$ws.Range("F29").Value = 174.04

# Row 30: Mackay
$ws.Range("A30").Value = 45952
$ws.Range("D30").Value = 167.68
$ws.Range("E30").Value = 163.43
$ws.Range("F30").Value = 173.43
$ws.Range("G30").Value = 163.26

# Row 31: Townsville
$ws.Range("A31").Value = 45952
$ws.Range("D31").Value = 166.4
$ws.Range("E31").Value = 164.97
$ws.Range("F31").Value = 174.97

# Row 35: Adelaide
$ws.Range("A35").Value = 45953
$ws.Range("D35").Value = 160.13
$ws.Range("E35").Value = 157.07
$ws.Range("F35").Value = 166.07

# Row 36: Adelaide
$ws.Range("A36").Value = 45952
$ws.Range("D36").Value = 160.52000000000001
$ws.Range("E36").Value = 157.33000000000001
$ws.Range("F36").Value = 166.33

# Row 40: Burnie
$ws.Range("A40").Value = 45953
$ws.Range("D40").Value = 165.78
$ws.Range("E40").Value = 162.54
$ws.Range("F40").Value = 172.54

# Row 41: Hobart
$ws.Range("A41").Value = 45953
$ws.Range("D41").Value = 165.5
$ws.Range("E41").Value = 162.96
$ws.Range("F41").Value = 172.96

# Row 42: Burnie
$ws.Range("A42").Value = 45952
$ws.Range("D42").Value = 166.17
$ws.Range("E42").Value = 162.81
$ws.Range("F42").Value = 172.81

# Row 43: Hobart
$ws.Range("A43").Value = 45952
$ws.Range("D43").Value = 165.89
$ws.Range("E43").Value = 163.22999999999999
$ws.Range("F43").Value = 173.23

# Row 47: Geelong
$ws.Range("A47").Value = 45953
$ws.Range("D47").Value = 160.19999999999999
$ws.Range("E47").Value = 158.82
$ws.Range("F47").Value = 168.82

# Row 48: Melbourne
$ws.Range("A48").Value = 45953
$ws.Range("D48").Value = 160.18
$ws.Range("E48").Value = 158.99
$ws.Range("F48").Value = 168.99

# Row 49: Geelong
$ws.Range("A49").Value = 45952
$ws.Range("D49").Value = 160.86000000000001
$ws.Range("E49").Value = 159.02000000000001
$ws.Range("F49").Value = 169.02

# Row 50: Melbourne
$ws.Range("A50").Value = 45952
$ws.Range("D50").Value = 160.84
$ws.Range("E50").Value = 159.19
$ws.Range("F50").Value = 169.19

# Row 54: Broome
$ws.Range("A54").Value = 45953
$ws.Range("D54").Value = 175.96
$ws.Range("E54").Value = 172.84
$ws.Range("F54").Value = 182.84

# Row 55: Esperance
$ws.Range("A55").Value = 45953
$ws.Range("D55").Value = 163.63
$ws.Range("E55").Value = 170.13
$ws.Range("F55").Value = 180.13

# Row 56: Geraldton
$ws.Range("A56").Value = 45953
$ws.Range("D56").Value = 165.91

# Row 57: Kalgoorlie
$ws.Range("A57").Value = 45953
$ws.Range("D57").Value = 165.58
$ws.Range("E57").Value = 164.4

# Row 58: Perth
$ws.Range("A58").Value = 45953
$ws.Range("D58").Value = 161.47999999999999
$ws.Range("E58").Value = 160.44999999999999
$ws.Range("F58").Value = 170.45

# Row 59: Port Hedland
$ws.Range("A59").Value = 45953
$ws.Range("D59").Value = 168.3
$ws.Range("E59").Value = 171.03

# Row 60: Broome
$ws.Range("A60").Value = 45952
$ws.Range("D60").Value = 176.36
$ws.Range("E60").Value = 173.11
$ws.Range("F60").Value = 183.11

# Row 61: Esperance
$ws.Range("A61").Value = 45952
$ws.Range("D61").Value = 164.01
$ws.Range("E61").Value = 170.5
$ws.Range("F61").Value = 180.5

# Row 62: Geraldton
$ws.Range("A62").Value = 45952
$ws.Range("D62").Value = 166.3

# Row 63: Kalgoorlie
$ws.Range("A63").Value = 45952
$ws.Range("D63").Value = 165.96
$ws.Range("E63").Value = 164.77

# Row 64: Perth
$ws.Range("A64").Value = 45952
$ws.Range("D64").Value = 161.87
$ws.Range("E64").Value = 160.82
$ws.Range("F64").Value = 170.82

# Row 65: Port Hedland
$ws.Range("A65").Value = 45952
$ws.Range("D65").Value = 168.69
$ws.Range("E65").Value = 171.29
